# Helper: PowerPoint's Shape.Left/Top (and other position/size) properties
# are expressed in points, but OOXML stores offsets in EMU (1 pt = 12700 EMU).
# The COM layer marshals the point value through a single-precision float
# before re-quantizing to EMU, so naive `emu / 12700.0` can land one EMU off
# after the round-trip. This helper searches (in 1e-6 pt steps) for a point
# value whose single-precision round-trip reproduces the exact target EMU.
function EmuToPt {
    param([double]$Emu)

    $base = $Emu / 12700.0
    for ($d = -500; $d -le 500; $d++) {
        $cand = $base + ($d * 0.000001)
        $f = [double]([single]$cand)
        $gotEmu = [math]::Floor($f * 12700)
        if ($gotEmu -eq $Emu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Existing shape: "Rectangle : coins arrondis 3" (id 4) — the "GW-P 2020" badge.
$orig = $s.Shapes.Item(1)

# Duplicate it first (before retexting/moving) so the copy starts out as an
# exact clone of the original "GW-P 2020" badge, matching the new shape
# introduced by the edit ("Rectangle : coins arrondis 1", id 2).
$copy = $orig.Duplicate()
$copy.Name = "Rectangle : coins arrondis 1"
$copy.Left = EmuToPt(4625196)
$copy.Top = EmuToPt(2402457)
$copy.TextFrame.TextRange.Text = "GW-P 2020C"

# Move/retext the original shape to its new spot and wording.
$orig.Left = EmuToPt(4478547)
$orig.Top = EmuToPt(560717)
$orig.TextFrame.TextRange.Text = "GW-P 2020C"
